$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entity types")

# Insert two new entity types ("PrintMedia" and "Quantity") alphabetically
# between "Person" (row 19) and "Sport" (row 20), pushing the existing
# rows 20-24 down to 22-26.
$ws.Rows.Item(20).Resize(2).Insert()

$ws.Range("A20").Value = "PrintMedia"
$ws.Range("A21").Value = "Quantity"

$ws.Range("A26").Select()
